$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# Fill in row 9 with data for the new "Implementar PilaCL" task
$ws.Range("A9").Value = "Implementar PilaCL"
$ws.Range("B9").Value = 25
$ws.Range("C9").Value = 22
$ws.Range("D9").Value = 0.006944444444444444
$ws.Range("E9").Value = 0.8333333333333334
$ws.Range("F9").Value = 0.8409722222222222
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

# Update the selected cell to A9
$ws.Range("A9").Select()

$wb.Save()
